$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend header row 1 with P1=14, Q1=15, copying the style from O1 (B1:O1 use style "1": bold + border)
$ws.Range("O1").Copy() | Out-Null
$ws.Range("P1:Q1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("P1").Value2 = 14
$ws.Range("Q1").Value2 = 15

# For each data row (2-25), swap columns I<->K and M<->O, then fill new columns P and Q with 2
for ($r = 2; $r -le 25; $r++) {
    $iVal = $ws.Cells.Item($r, 9).Value2   # column I
    $kVal = $ws.Cells.Item($r, 11).Value2  # column K
    $mVal = $ws.Cells.Item($r, 13).Value2  # column M
    $oVal = $ws.Cells.Item($r, 15).Value2  # column O

    $ws.Cells.Item($r, 9).Value2 = $kVal   # I = old K
    $ws.Cells.Item($r, 11).Value2 = $iVal  # K = old I
    $ws.Cells.Item($r, 13).Value2 = $oVal  # M = old O
    $ws.Cells.Item($r, 15).Value2 = $mVal  # O = old M

    $ws.Cells.Item($r, 16).Value2 = 2      # P
    $ws.Cells.Item($r, 17).Value2 = 2      # Q
}
